# Actualizacion automatica del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill resultado/profit for two previously-pending rows ---
$ws.Cells.Item(171, 7).Value = "Fallo"
$ws.Cells.Item(171, 8).Value = -1

$ws.Cells.Item(173, 7).Value = "Fallo"
$ws.Cells.Item(173, 8).Value = -1

# --- Append new pending matches (rows 179-187) ---
# Each new row keeps resultado/profit empty (pending), same as source rows
# before they are settled.

$row = 179
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14851669
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Jan-Lennard Struff"
$ws.Cells.Item($row, 4).Value = "Corentin Moutet"
$ws.Cells.Item($row, 5).Value = "Gana Corentin Moutet"
$ws.Cells.Item($row, 6).Value = 1.73

$row = 180
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14852137
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Daniel Dutra Da Silva"
$ws.Cells.Item($row, 4).Value = "Joao Eduardo Schiessl"
$ws.Cells.Item($row, 5).Value = "Gana Joao Eduardo Schiessl"
$ws.Cells.Item($row, 6).Value = 2

$row = 181
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14862976
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Alice Rame"
$ws.Cells.Item($row, 4).Value = "Martina Colmegna"
$ws.Cells.Item($row, 5).Value = "Gana Martina Colmegna"
$ws.Cells.Item($row, 6).Value = 4

$row = 182
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14900104
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Vladyslav Orlov"
$ws.Cells.Item($row, 4).Value = "Abel Forger"
$ws.Cells.Item($row, 5).Value = "Gana Abel Forger"
$ws.Cells.Item($row, 6).Value = 2

$row = 183
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14896439
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Yassine Dlimi"
$ws.Cells.Item($row, 4).Value = "Florent Bax"
$ws.Cells.Item($row, 5).Value = "Gana Yassine Dlimi"
$ws.Cells.Item($row, 6).Value = 2.2

$row = 184
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14899184
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Eric Vanshelboim"
$ws.Cells.Item($row, 4).Value = "Tiago Pereira"
$ws.Cells.Item($row, 5).Value = "Gana Eric Vanshelboim"
$ws.Cells.Item($row, 6).Value = 3.4

$row = 185
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14899531
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Luca Wiedenmann"
$ws.Cells.Item($row, 4).Value = "Etienne Donnet"
$ws.Cells.Item($row, 5).Value = "Gana Luca Wiedenmann"
$ws.Cells.Item($row, 6).Value = 2.62

$row = 186
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14899183
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Mathys Erhard"
$ws.Cells.Item($row, 4).Value = "Jack Loge"
$ws.Cells.Item($row, 5).Value = "Gana Jack Loge"
$ws.Cells.Item($row, 6).Value = 2.25

$row = 187
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = 14899199
$ws.Cells.Item($row, 2).Value = "2025-10-17"
$ws.Cells.Item($row, 3).Value = "Gianmarco Ferrari"
$ws.Cells.Item($row, 4).Value = "Gianluca Cadenasso"
$ws.Cells.Item($row, 5).Value = "Gana Gianmarco Ferrari"
$ws.Cells.Item($row, 6).Value = 2.75
